$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Hunk 0: ALC!row6 (Days of Chunder)
$ws_ALC.Range("H6").Value = 44609.4
$ws_ALC.Range("I6").Value = 331.375
$ws_ALC.Range("K6").Value = 994.125
$ws_ALC.Range("M6").Value = -882.125

# Hunk 1: ALC!row137 (Cutting Edge of Culinary Quality)
$ws_ALC.Range("H137").Value = 1450510
$ws_ALC.Range("I137").Value = 1923951.2
$ws_ALC.Range("J137").Value = 2336.7058
$ws_ALC.Range("K137").Value = 5771853.6
$ws_ALC.Range("L137").Value = 7010.117400000001
$ws_ALC.Range("M137").Value = -5769303.6
$ws_ALC.Range("N137").Value = -12110.1174

# Hunk 2: ALC!row138 (All-night Crafting)
$ws_ALC.Range("H138").Value = 2316418.5
$ws_ALC.Range("I138").Value = 1169.8276
$ws_ALC.Range("J138").Value = 11908164
$ws_ALC.Range("K138").Value = 3509.4828
$ws_ALC.Range("L138").Value = 35724492
$ws_ALC.Range("M138").Value = 1630.5172
$ws_ALC.Range("N138").Value = -35734772

# Hunk 3: ARM!row2 (Ain't Got No Ingots)
$ws_ARM.Range("H2").Value = 4101.3
$ws_ARM.Range("I2").Value = 4400
$ws_ARM.Range("J2").Value = 4026.625
$ws_ARM.Range("K2").Value = 4400
$ws_ARM.Range("L2").Value = 4026.625
$ws_ARM.Range("M2").Value = -4287
$ws_ARM.Range("N2").Value = -4252.625

# Hunk 4: ARM!row32 (Ingot We Trust)
$ws_ARM.Range("H32").Value = 2911.81
$ws_ARM.Range("I32").Value = 2507.8052
$ws_ARM.Range("J32").Value = 4264.3477
$ws_ARM.Range("K32").Value = 2507.8052
$ws_ARM.Range("L32").Value = 4264.3477
$ws_ARM.Range("M32").Value = -2220.8052
$ws_ARM.Range("N32").Value = -4838.3477

# Hunk 5: ARM!row44 (Very Slow Array)
$ws_ARM.Range("H44").Value = 16666.666
$ws_ARM.Range("I44").Value = 10000
$ws_ARM.Range("K44").Value = 10000
$ws_ARM.Range("M44").Value = -9512

# Hunk 6: ARM!row45 (Hollow Hallmarks)
$ws_ARM.Range("H45").Value = 1693.75
$ws_ARM.Range("I45").Value = 1090
$ws_ARM.Range("J45").Value = 2700
$ws_ARM.Range("K45").Value = 1090
$ws_ARM.Range("L45").Value = 2700
$ws_ARM.Range("M45").Value = -713
$ws_ARM.Range("N45").Value = -3454

# Hunk 7: ARM!row61 (Dealing with the Tough Stuff)
$ws_ARM.Range("H61").Value = 18556850
$ws_ARM.Range("I61").Value = 21762008
$ws_ARM.Range("J61").Value = 127189.125
$ws_ARM.Range("K61").Value = 21762008
$ws_ARM.Range("L61").Value = 127189.125
$ws_ARM.Range("M61").Value = -21761796
$ws_ARM.Range("N61").Value = -127613.125

# Hunk 8: ARM!row74 (As the Bolt Flies)
$ws_ARM.Range("H74").Value = 6844732.5
$ws_ARM.Range("I74").Value = 8573617
$ws_ARM.Range("K74").Value = 8573617
$ws_ARM.Range("M74").Value = -8572743

# Hunk 9: ARM!row77 (Heavy Metal Banned (L))
$ws_ARM.Range("H77").Value = 6844732.5
$ws_ARM.Range("I77").Value = 8573617
$ws_ARM.Range("K77").Value = 42868085
$ws_ARM.Range("M77").Value = -42863717

# Hunk 10: ARM!row110 (Scheduled Maintenance)
$ws_ARM.Range("H110").Value = 1283.7142
$ws_ARM.Range("I110").Value = 1157.2
$ws_ARM.Range("K110").Value = 1157.2
$ws_ARM.Range("M110").Value = 887.8

# Hunk 11: ARM!row116 (No Scope)
$ws_ARM.Range("H116").Value = 4101.3
$ws_ARM.Range("I116").Value = 4400
$ws_ARM.Range("J116").Value = 4026.625
$ws_ARM.Range("K116").Value = 4400
$ws_ARM.Range("L116").Value = 4026.625
$ws_ARM.Range("M116").Value = -2106
$ws_ARM.Range("N116").Value = -8614.625

# Hunk 12: ARM!row136 (Metal with Mettle)
$ws_ARM.Range("H136").Value = 18556850
$ws_ARM.Range("I136").Value = 21762008
$ws_ARM.Range("J136").Value = 127189.125
$ws_ARM.Range("K136").Value = 65286024
$ws_ARM.Range("L136").Value = 381567.375
$ws_ARM.Range("M136").Value = -65283474
$ws_ARM.Range("N136").Value = -386667.375

# Hunk 13: BSM!row3 (Hells Bells)
$ws_BSM.Range("H3").Value = 4101.3
$ws_BSM.Range("I3").Value = 4400
$ws_BSM.Range("J3").Value = 4026.625
$ws_BSM.Range("K3").Value = 4400
$ws_BSM.Range("L3").Value = 4026.625
$ws_BSM.Range("M3").Value = -4286
$ws_BSM.Range("N3").Value = -4254.625

# Hunk 14: BSM!row107 (The Gold Experience)
$ws_BSM.Range("H107").Value = 2035.6666
$ws_BSM.Range("I107").Value = 1859.2142
$ws_BSM.Range("J107").Value = 2653.25
$ws_BSM.Range("K107").Value = 1859.2142
$ws_BSM.Range("L107").Value = 2653.25
$ws_BSM.Range("M107").Value = 60.78580000000011
$ws_BSM.Range("N107").Value = -6493.25

# Hunk 15: BSM!row134 (Ruthenium Supremium)
$ws_BSM.Range("H134").Value = 1360.1904
$ws_BSM.Range("I134").Value = 850.3333
$ws_BSM.Range("K134").Value = 2550.9999
$ws_BSM.Range("M134").Value = -15.9998999999998

# Hunk 16: CRP!row31 (Wall Not Found)
$ws_CRP.Range("H31").Value = 4096.926
$ws_CRP.Range("I31").Value = 1663.7222
$ws_CRP.Range("J31").Value = 8963.333000000001
$ws_CRP.Range("K31").Value = 1663.7222
$ws_CRP.Range("L31").Value = 8963.333000000001
$ws_CRP.Range("M31").Value = -1368.7222
$ws_CRP.Range("N31").Value = -9553.333000000001

# Hunk 17: CRP!row34 (Armoires of the Rich and Famous)
$ws_CRP.Range("H34").Value = 4096.926
$ws_CRP.Range("I34").Value = 1663.7222
$ws_CRP.Range("J34").Value = 8963.333000000001
$ws_CRP.Range("K34").Value = 1663.7222
$ws_CRP.Range("L34").Value = 8963.333000000001
$ws_CRP.Range("M34").Value = -1461.7222
$ws_CRP.Range("N34").Value = -9367.333000000001

# Hunk 18: CRP!row50 (The Arsenal of Theocracy)
$ws_CRP.Range("H50").Value = 21909.715
$ws_CRP.Range("J50").Value = 21909.715
$ws_CRP.Range("L50").Value = 21909.715
$ws_CRP.Range("N50").Value = -23159.715

# Hunk 19: CRP!row59 (Bow Down to Magic)
$ws_CRP.Range("H59").Value = 33111.125
$ws_CRP.Range("J59").Value = 33111.125
$ws_CRP.Range("L59").Value = 33111.125
$ws_CRP.Range("N59").Value = -35401.125

# Hunk 20: CRP!row60 (Bowing to Greater Power)
$ws_CRP.Range("H60").Value = 6500
$ws_CRP.Range("I60").Value = 3666.6667
$ws_CRP.Range("K60").Value = 3666.6667
$ws_CRP.Range("M60").Value = -3155.6667

# Hunk 21: CRP!row99 (O Pine)
$ws_CRP.Range("H99").Value = 7080
$ws_CRP.Range("I99").Value = 1600
$ws_CRP.Range("J99").Value = 10733.333
$ws_CRP.Range("K99").Value = 1600
$ws_CRP.Range("L99").Value = 10733.333
$ws_CRP.Range("M99").Value = -102
$ws_CRP.Range("N99").Value = -13729.333

# Hunk 22: CRP!row107 (Built to Last)
$ws_CRP.Range("H107").Value = 313.13635
$ws_CRP.Range("I107").Value = 323.2857
$ws_CRP.Range("J107").Value = 100
$ws_CRP.Range("K107").Value = 323.2857
$ws_CRP.Range("L107").Value = 100
$ws_CRP.Range("M107").Value = 1596.7143
$ws_CRP.Range("N107").Value = -3940

# Hunk 23: CRP!row126 (A Better Conductor)
$ws_CRP.Range("H126").Value = 7080
$ws_CRP.Range("I126").Value = 1600
$ws_CRP.Range("J126").Value = 10733.333
$ws_CRP.Range("K126").Value = 4800
$ws_CRP.Range("L126").Value = 32199.999
$ws_CRP.Range("M126").Value = -2330
$ws_CRP.Range("N126").Value = -37139.999

# Hunk 24: CRP!row132 (Hull Lotta Damage)
$ws_CRP.Range("H132").Value = 40700
$ws_CRP.Range("I132").Value = 2180.1052
$ws_CRP.Range("J132").Value = 145254
$ws_CRP.Range("K132").Value = 6540.3156
$ws_CRP.Range("L132").Value = 435762
$ws_CRP.Range("M132").Value = -4010.3156
$ws_CRP.Range("N132").Value = -440822

# Hunk 25: CRP!row134 (Wood You Be Quiet)
$ws_CRP.Range("H134").Value = 31596.695
$ws_CRP.Range("I134").Value = 1893.7894
$ws_CRP.Range("J134").Value = 64794.06
$ws_CRP.Range("K134").Value = 5681.3682
$ws_CRP.Range("L134").Value = 194382.18
$ws_CRP.Range("M134").Value = -3146.3682
$ws_CRP.Range("N134").Value = -199452.18

# Hunk 26: CUL!row29 (For Crumbs' Sake)
$ws_CUL.Range("H29").Value = 370
$ws_CUL.Range("I29").Value = 26
$ws_CUL.Range("J29").Value = 800
$ws_CUL.Range("K29").Value = 78
$ws_CUL.Range("L29").Value = 2400
$ws_CUL.Range("M29").Value = 199
$ws_CUL.Range("N29").Value = -2954

# Hunk 27: CUL!row34 (Fever Pitch)
$ws_CUL.Range("H34").Value = 1007.7857
$ws_CUL.Range("I34").Value = 380
$ws_CUL.Range("J34").Value = 1356.5555
$ws_CUL.Range("K34").Value = 1140
$ws_CUL.Range("L34").Value = 4069.6665
$ws_CUL.Range("M34").Value = -1056
$ws_CUL.Range("N34").Value = -4237.666499999999

# Hunk 28: CUL!row92 (Oh No Udon)
$ws_CUL.Range("H92").Value = 858.7059
$ws_CUL.Range("I92").Value = 1010
$ws_CUL.Range("J92").Value = 795.6667
$ws_CUL.Range("K92").Value = 3030
$ws_CUL.Range("L92").Value = 2387.0001
$ws_CUL.Range("M92").Value = -1782
$ws_CUL.Range("N92").Value = -4883.0001

# Hunk 29: CUL!row124 (Bobbing for Compliments)
$ws_CUL.Range("H124").Value = 973.55
$ws_CUL.Range("I124").Value = 490
$ws_CUL.Range("K124").Value = 1470
$ws_CUL.Range("M124").Value = 3440

# Hunk 30: CUL!row131 (The Mountain Steeped)
$ws_CUL.Range("H131").Value = 1015.5469
$ws_CUL.Range("J131").Value = 1101.1608
$ws_CUL.Range("L131").Value = 3303.4824
$ws_CUL.Range("N131").Value = -13383.4824

# Hunk 31: LTW!row3 (Underneath It All)
$ws_LTW.Range("H3").Value = 500000
$ws_LTW.Range("J3").Value = 0
$ws_LTW.Range("L3").Value = 0
$ws_LTW.Range("N3").ClearContents()

# Hunk 32: LTW!row15 (The Bards' Guards)
$ws_LTW.Range("H15").Value = 500000
$ws_LTW.Range("J15").Value = 0
$ws_LTW.Range("L15").Value = 0
$ws_LTW.Range("N15").ClearContents()

# Hunk 33: LTW!row100 (Tiger in the Sack)
$ws_LTW.Range("H100").Value = 1693.6538
$ws_LTW.Range("I100").Value = 1613.5
$ws_LTW.Range("K100").Value = 1613.5
$ws_LTW.Range("M100").Value = -1072.5

# Hunk 34: WVR!row96 (Skills on Display)
$ws_WVR.Range("H96").Value = 1500
$ws_WVR.Range("I96").Value = 0
$ws_WVR.Range("J96").Value = 1500
$ws_WVR.Range("K96").Value = 0
$ws_WVR.Range("L96").Value = 1500
$ws_WVR.Range("M96").ClearContents()
$ws_WVR.Range("N96").Value = -4246

# Hunk 35: WVR!row136 (Weaving the Envelope)
$ws_WVR.Range("H136").Value = 23330.252
$ws_WVR.Range("I136").Value = 16773.709
$ws_WVR.Range("J136").Value = 37347.69
$ws_WVR.Range("K136").Value = 50321.12699999999
$ws_WVR.Range("L136").Value = 112043.07
$ws_WVR.Range("M136").Value = -47771.12699999999
$ws_WVR.Range("N136").Value = -117143.07
